$d = $word.ActiveDocument

function Get-ParagraphByStart($doc, $marker) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t.StartsWith($marker)) {
            return $p.Range
        }
    }
    return $null
}

# Change 1: merge "Backlogs are " + "incomplete...sprint 1 " into a single run
$marker1 = 'Backlogs are '
$para1 = Get-ParagraphByStart $d $marker1
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="7477486B" w14:textId="3AD1E9F1" w:rsidR="00EA0347" w:rsidRDefault="00EA0347" w:rsidP="00EA0347"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">Backlogs are incomplete—System Product Backlog is not up to date—only sprint 1 </w:t></w:r><w:r w:rsidR="00AA03A0"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t>Assignments are listed—Sprint 4 backlog is incomplete—presentation in fact only showed sprint 3 backlog, not sprint 4!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para1.InsertXML($xml1)

# Change 2: split "Sprint 3 Team Evaluation..." into "Sprint " + "4" + " Team Evaluation..."
$marker2 = 'Sprint 3 Team Evaluation is very good. Have you gotten the DB connected?'
$para2 = Get-ParagraphByStart $d $marker2
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0FC2AF90" w14:textId="77777777" w:rsidR="00EA0347" w:rsidRDefault="00EA0347" w:rsidP="00EA0347"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve">Sprint </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t xml:space="preserve"> Team Evaluation is very good. Have you gotten the DB connected?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para2.InsertXML($xml2)

# Change 3: merge the 4 "Project plan needs..." runs into a single run
$marker3 = 'Project plan needs to be updated to show all '
$para3 = Get-ParagraphByStart $d $marker3
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3AD5E4AE" w14:textId="122ED2D2" w:rsidR="00EA0347" w:rsidRDefault="00EA0347" w:rsidP="00EA0347"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr><w:t>Project plan needs to be updated to show all dates on the timeline—if you need help doing that let me know—Sprint 4 items still say 0% complete</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para3.InsertXML($xml3)

Write-Host "Done applying edits"
